$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 132, shifting existing rows 132:237 down to 133:238
$ws.Rows(132).Insert()

# Populate the newly inserted row 132 with the new record
$ws.Range("A132").Value = 8
$ws.Range("B132").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C132").Value = 'Coquimbo'
$ws.Range("D132").Value = 44729
$ws.Range("E132").Value = 4
$ws.Range("F132").Value = 100112031
$ws.Range("G132").Value = 'Poroto verde'
$ws.Range("H132").Value = 'Magnum'
$ws.Range("I132").Value = 'Primera'
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 25000
$ws.Range("L132").Value = 26000
$ws.Range("M132").Value = 25500
$ws.Range("N132").Value = '$/malla 25 kilos'
$ws.Range("O132").Value = 'Perú'
$ws.Range("P132").Value = 1020
$ws.Range("Q132").Value = 25
$ws.Range("R132").Value = 'Hortaliza'
